$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.645.29"
$ws.Range("E2").Value = "  +2.72%  "

Set-TextValue $ws.Range("D3") "1.912.26"
$ws.Range("E3").Value = "  +5.56%  "

Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "313.95"
$ws.Range("E5").Value = "  +1.48%  "

Set-TextValue $ws.Range("D6") "0.9996"
$ws.Range("E6").Value = "  -0.08%  "

Set-TextValue $ws.Range("D7") "0.5047"
$ws.Range("E7").Value = "  +2.18%  "

Set-TextValue $ws.Range("D8") "0.3970"
$ws.Range("E8").Value = "  +2.19%  "

Set-TextValue $ws.Range("D9") "0.09656"
$ws.Range("E9").Value = "  -1.03%  "

Set-TextValue $ws.Range("D10") "1.162"
$ws.Range("E10").Value = "  +5.57%  "

Set-TextValue $ws.Range("D11") "41.79"
$ws.Range("E11").Value = "  +2.33%  "

Set-TextValue $ws.Range("D12") "6.551"
$ws.Range("E12").Value = "  +1.91%  "

Set-TextValue $ws.Range("D13") "21.24"
$ws.Range("E13").Value = "  +3.73%  "

Set-TextValue $ws.Range("D14") "1.918.91"
$ws.Range("E14").Value = "  +6.15%  "

Set-TextValue $ws.Range("D15") "7.565"
$ws.Range("E15").Value = "  +3.73%  "

Set-TextValue $ws.Range("D16") "1.000"
$ws.Range("E16").Value = "  -0.01%  "

Set-TextValue $ws.Range("D17") "0.00001137"
$ws.Range("E17").Value = "  +0.35%  "

Set-TextValue $ws.Range("D18") "94.17"
$ws.Range("E18").Value = "  +1.73%  "

Set-TextValue $ws.Range("D19") "0.06638"
$ws.Range("E19").Value = "  +0.57%  "

Set-TextValue $ws.Range("D20") "18.02"
$ws.Range("E20").Value = "  +5.55%  "

$ws.Range("E22").Value = "  +5.90%  "

Set-TextValue $ws.Range("D23") "28.696.63"
$ws.Range("E23").Value = "  +2.71%  "

Set-TextValue $ws.Range("D24") "11.45"
$ws.Range("E24").Value = "  +2.69%  "

Set-TextValue $ws.Range("D25") "2.284"
$ws.Range("E25").Value = "  +1.80%  "

Set-TextValue $ws.Range("D26") "2.772"
$ws.Range("E26").Value = "  +15.81%  "

Set-TextValue $ws.Range("D27") "2.128.62"
$ws.Range("E27").Value = "  +5.60%  "

Set-TextValue $ws.Range("D28") "21.51"
$ws.Range("E28").Value = "  +4.33%  "

Set-TextValue $ws.Range("D29") "159.65"
$ws.Range("E29").Value = "  +1.34%  "

Set-TextValue $ws.Range("D30") "128.76"
$ws.Range("E30").Value = "  +1.28%  "

Set-TextValue $ws.Range("D31") "1.110"
$ws.Range("E31").Value = "  +6.68%  "

Set-TextValue $ws.Range("D32") "0.1076"
$ws.Range("E32").Value = "  +1.52%  "

Set-TextValue $ws.Range("D33") "5.744"
$ws.Range("E33").Value = "  +2.96%  "

Set-TextValue $ws.Range("D34") "3.616"
$ws.Range("E34").Value = "  -0.45%  "

Set-TextValue $ws.Range("D35") "9.799"
$ws.Range("E35").Value = "  +8.67%  "

Set-TextValue $ws.Range("D36") "0.06817"
$ws.Range("E36").Value = "  +0.95%  "

Set-TextValue $ws.Range("D37") "0.02445"
$ws.Range("E37").Value = "  +5.33%  "

Set-TextValue $ws.Range("D38") "0.2219"
$ws.Range("E38").Value = "  +4.17%  "

Set-TextValue $ws.Range("D39") "5.096"
$ws.Range("E39").Value = "  +3.18%  "

Set-TextValue $ws.Range("D40") "11.65"
$ws.Range("E40").Value = "  +3.55%  "

Set-TextValue $ws.Range("D41") "0.6423"
$ws.Range("E41").Value = "  +3.79%  "

Set-TextValue $ws.Range("D42") "1.199"
$ws.Range("E42").Value = "  +4.75%  "

Set-TextValue $ws.Range("D43") "0.9988"
$ws.Range("E43").Value = "  -0.11%  "

Set-TextValue $ws.Range("D44") "13.86"
$ws.Range("E44").Value = "  +5.90%  "

Set-TextValue $ws.Range("D45") "0.6090"
$ws.Range("E45").Value = "  +3.82%  "

Set-TextValue $ws.Range("D46") "1.287"
$ws.Range("E46").Value = "  +0.40%  "

Set-TextValue $ws.Range("D47") "3.664"
$ws.Range("E47").Value = "  -0.77%  "

Set-TextValue $ws.Range("D48") "2.044"
$ws.Range("E48").Value = "  +5.65%  "

Set-TextValue $ws.Range("D49") "124.83"
$ws.Range("E49").Value = "  +1.94%  "

Set-TextValue $ws.Range("D50") "1.213"
$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "78.45"
$ws.Range("E51").Value = "  +6.80%  "
